{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph immediately preceding it) that followed the\n// last bibliography entry (\"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear. S\u00e3o Paulo:\n// Thomson, 2007.\"), while leaving the blank paragraph / page-break\n// paragraph that come after the footer untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the bibliography paragraph that ends this run of text.\nconst anchorText = \"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear. S\u00e3o Paulo: Thomson, 2007.\";\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find anchor paragraph for the bibliography entry.\");\n}\n\nconst jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst copyrightText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\n// Right after the anchor paragraph we expect: blank paragraph, the\n// \"Ver no Jupiter...\" paragraph, then the \"\u00a9 2020 ...\" paragraph.\nconst blankIndex = anchorIndex + 1;\nconst jupiterIndex = anchorIndex + 2;\nconst copyrightIndex = anchorIndex + 3;\n\nif (\n  blankIndex >= items.length ||\n  jupiterIndex >= items.length ||\n  copyrightIndex >= items.length ||\n  items[blankIndex].text !== \"\" ||\n  items[jupiterIndex].text !== jupiterText ||\n  items[copyrightIndex].text !== copyrightText\n) {\n  throw new Error(\"Document structure did not match the expected footer block.\");\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[jupiterIndex].delete();\nitems[blankIndex].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph immediately preceding it) that followed the\n# last bibliography entry (\"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear. S\u00e3o Paulo:\n# Thomson, 2007.\"), while leaving the blank paragraph / page-break\n# paragraph that come after the footer untouched.\n\n$d = $word.ActiveDocument\n\n$anchorText    = \"uma introdu\u00e7\u00e3o \u00e0 \u00e1lgebra linear. S\u00e3o Paulo: Thomson, 2007.\" + [char]13\n$jupiterText   = \"Ver no Jupiter Salvar em pdf Salvar em docx\" + [char]13\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\" + [char]13\n\n# Locate the bibliography paragraph that anchors the footer block.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text -eq $anchorText) {\n    $anchorIndex = $i\n    break\n  }\n}\n\nif ($anchorIndex -eq -1) {\n  throw \"Could not find anchor paragraph for the bibliography entry.\"\n}\n\n$blankIndex     = $anchorIndex + 1\n$jupiterIndex   = $anchorIndex + 2\n$copyrightIndex = $anchorIndex + 3\n\nif ($d.Paragraphs.Item($blankIndex).Range.Text -ne ([char]13) `\n    -or $d.Paragraphs.Item($jupiterIndex).Range.Text -ne $jupiterText `\n    -or $d.Paragraphs.Item($copyrightIndex).Range.Text -ne $copyrightText) {\n  throw \"Document structure did not match the expected footer block.\"\n}\n\n# Delete from the end backwards so earlier indices stay valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($jupiterIndex).Range.Delete()\n$d.Paragraphs.Item($blankIndex).Range.Delete()\n"}
